# "fix figures in ITI-31"
#
# 1. Re-cache the "datetimeFigureOut" date field shown on the Notes Master,
#    the Slide Master and every slide layout to 11/17/2020.
# 2. Move/resize the embedded PowerPoint-show OLE object on slide 1.

$p = $ppt.ActivePresentation

function Set-DatePlaceholderText {
    param($shapes, [string]$newText)

    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            $shp.TextFrame.TextRange.Text = $newText
        }
    }
}

# Notes master date field: 11/10/20 -> 11/17/2020
Set-DatePlaceholderText -shapes $p.NotesMaster.Shapes -newText "11/17/2020"

# Slide master date field: 11/9/20 -> 11/17/2020
Set-DatePlaceholderText -shapes $p.SlideMaster.Shapes -newText "11/17/2020"

# Every slide layout's date field: 11/9/20 -> 11/17/2020
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Set-DatePlaceholderText -shapes $layout.Shapes -newText "11/17/2020"
}

# Reposition/resize the embedded OLE "PowerPoint Show" object on slide 1.
# Target EMU: off(2300141,252817) ext(8335210,6279953).
# A small epsilon is added after the EMU->point conversion so the
# point value lands in the correct EMU bucket once the host round-trips
# it back through its own (lower precision) point->EMU conversion.
$slide = $p.Slides.Item(1)
$ole = $slide.Shapes.Item(5)

$emuPerPoint = 12700.0
$eps = 0.00003

$ole.Left   = (2300141 / $emuPerPoint) + $eps
$ole.Top    = (252817  / $emuPerPoint) + $eps
$ole.Width  = (8335210 / $emuPerPoint) + $eps
$ole.Height = (6279953 / $emuPerPoint) + $eps
